$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: column width math.
# This COM-interop's ColumnWidth setter snaps the stored OOXML <col width>
# value to the nearest 1/6th of a character and offsets it by 5/6 (the
# standard Excel cell-padding constant). To land as close as possible on a
# desired final OOXML width, solve for the ColumnWidth input that reproduces
# it: input = desired_xml_width - 5/6
# ---------------------------------------------------------------------------

$wOverviewStatus = 29.9777047293527 - (5/6)   # Overview!E:F  (status text got longer)
$wTableStatus    = 29.9777047293527 - (5/6)   # zh-cn/de-de Status column (C)
$wTableFile      = 40 - (5/6)                 # zh-cn/de-de Latest Target/Handback File columns (I,J)

# ---------------------------------------------------------------------------
# Overview sheet: handback status text changed for both language rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E1:F1").ColumnWidth = $wOverviewStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: populate handback columns (I/J/K) for the localized file
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I2").Value = "7c984b00-1d7f-429b-8bb7-65de61b23310.md"
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/7c984b00-1d7f-429b-8bb7-65de61b23310.md", [Type]::Missing, [Type]::Missing, "7c984b00-1d7f-429b-8bb7-65de61b23310.md") | Out-Null

$wsZhCn.Range("J2").Value = "7c984b00-1d7f-429b-8bb7-65de61b23310.2773e68546cbdfc602e930e1f106e4466926cca6.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-23 19:02:51"

$wsZhCn.Range("C1").ColumnWidth = $wTableStatus
$wsZhCn.Range("I1:J1").ColumnWidth = $wTableFile

# ---------------------------------------------------------------------------
# de-de sheet: populate handback columns (I/J/K) for the localized file
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I2").Value = "7c984b00-1d7f-429b-8bb7-65de61b23310.md"
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0faa8f81a7687ac0b7efe767efa8dd196af02340/e2e/7c984b00-1d7f-429b-8bb7-65de61b23310.md", [Type]::Missing, [Type]::Missing, "7c984b00-1d7f-429b-8bb7-65de61b23310.md") | Out-Null

$wsDeDe.Range("J2").Value = "7c984b00-1d7f-429b-8bb7-65de61b23310.2773e68546cbdfc602e930e1f106e4466926cca6.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-23 19:02:59"

$wsDeDe.Range("C1").ColumnWidth = $wTableStatus
$wsDeDe.Range("I1:J1").ColumnWidth = $wTableFile

$wb.Save()
